# Added if same name exists logic.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The single-letter labels in column A (rows 2-8) are being expanded to
# repeated-letter strings, simulating a "de-duplicate by appending the
# letter again when the same short name already exists" style fix.
$ws.Range("A2").Value = "AAAA"
$ws.Range("A3").Value = "BBBB"
$ws.Range("A4").Value = "CCCC"
$ws.Range("A5").Value = "DDDD"
$ws.Range("A6").Value = "EEEEE"
$ws.Range("A7").Value = "FFFFF"
$ws.Range("A8").Value = "GGGG"

# Move the active selection from A9 to A5.
$ws.Range("A5").Select()
